$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and Report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  34"
$ws.Range("C9").Value = "Report Covering the Week  8/21/2023  Through  8/27/2023"

# --- Cells changing from placeholder text to numeric: copy number style first ---
$ws.Range("G15").Copy($ws.Range("D27"))
$ws.Range("H15").Copy($ws.Range("E27"))

# --- Cells changing from numeric to placeholder text: copy placeholder style+value ---
$ws.Range("C15").Copy($ws.Range("G30"))
$ws.Range("E15").Copy($ws.Range("H30"))

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("G15").Value = 5
$ws.Range("L15").Value = 125
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 22
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = 4.761904761904
$ws.Range("I16").Value = 176
$ws.Range("J16").Value = 132
$ws.Range("K16").Value = 33.333333333333
$ws.Range("L16").Value = 131.578947368421
$ws.Range("M16").Value = 43.089430894308
$ws.Range("N16").Value = -74.231332357247
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 152
$ws.Range("J17").Value = 151
$ws.Range("K17").Value = 0.662251655629
$ws.Range("L17").Value = 14.285714285714
$ws.Range("M17").Value = 85.365853658536
$ws.Range("N17").Value = -25.490196078431
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 400
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 38.461538461538
$ws.Range("I18").Value = 133
$ws.Range("J18").Value = 113
$ws.Range("K18").Value = 17.699115044247
$ws.Range("L18").Value = 15.652173913043
$ws.Range("M18").Value = -16.352201257861
$ws.Range("N18").Value = -86.174636174636
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 51
$ws.Range("G19").Value = 64
$ws.Range("H19").Value = -20.3125
$ws.Range("I19").Value = 487
$ws.Range("J19").Value = 451
$ws.Range("K19").Value = 7.982261640798
$ws.Range("L19").Value = 79.044117647058
$ws.Range("M19").Value = 63.422818791946
$ws.Range("N19").Value = -18.288590604026
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 11
$ws.Range("E20").Value = -18.181818181818
$ws.Range("F20").Value = 37
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = 68.181818181818
$ws.Range("I20").Value = 175
$ws.Range("J20").Value = 132
$ws.Range("K20").Value = 32.575757575757
$ws.Range("L20").Value = 47.058823529411
$ws.Range("M20").Value = 22.377622377622
$ws.Range("N20").Value = -86.486486486486
$ws.Range("C21").Value = 43
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = 13.157894736842
$ws.Range("F21").Value = 142
$ws.Range("G21").Value = 146
$ws.Range("H21").Value = -2.739726027397
$ws.Range("I21").Value = 1141
$ws.Range("J21").Value = 994
$ws.Range("K21").Value = 14.788732394366
$ws.Range("L21").Value = 57.379310344827
$ws.Range("M21").Value = 39.657282741738
$ws.Range("N21").Value = -69.613848202396
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 6
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 11
$ws.Range("H22").Value = -54.545454545454
$ws.Range("I22").Value = 59
$ws.Range("J22").Value = 54
$ws.Range("K22").Value = 9.259259259259
$ws.Range("L22").Value = 136
$ws.Range("M22").Value = 96.666666666666
$ws.Range("C24").Value = 40
$ws.Range("D24").Value = 46
$ws.Range("E24").Value = -13.043478260869
$ws.Range("F24").Value = 199
$ws.Range("G24").Value = 148
$ws.Range("H24").Value = 34.459459459459
$ws.Range("I24").Value = 1302
$ws.Range("J24").Value = 931
$ws.Range("K24").Value = 39.849624060150
$ws.Range("L24").Value = 64.186633039092
$ws.Range("M24").Value = 110.339256865913
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -14.285714285714
$ws.Range("F25").Value = 39
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = 14.705882352941
$ws.Range("I25").Value = 334
$ws.Range("J25").Value = 376
$ws.Range("K25").Value = -11.170212765957
$ws.Range("L25").Value = 10.231023102310
$ws.Range("M25").Value = 0.602409638554
$ws.Range("G26").Value = 5
$ws.Range("L26").Value = 53.333333333333
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -57.142857142857
$ws.Range("J27").Value = 60
$ws.Range("K27").Value = -3.333333333333
$ws.Range("L27").Value = 38.095238095238

# --- Finalize values for cells that changed type (after style copy) ---
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -100
